# Suppression de lignes du dictionnaire de données
# Remove the last two rows of the data dictionary table:
#   - "Rôles préférés"
#   - "Liste films préférés"

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Delete from the bottom up so row indices stay valid.
$t.Rows.Item($t.Rows.Count).Delete()
$t.Rows.Item($t.Rows.Count).Delete()
